# Auto-generated Excel COM-interop script to apply the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.963.16"
$ws.Range("E2").Value = "  +4.75%  "
$ws.Range("D3").Value = "3.079.72"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.60"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.24"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.070.62"
$ws.Range("E8").Value = "  +2.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  +4.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.57"
$ws.Range("E11").Value = "  +7.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  +2.64%  "
$ws.Range("E13").Value = "  +4.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.43"
$ws.Range("E14").Value = "  +5.00%  "
$ws.Range("D16").Value = "3.589.55"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "3.078.25"
$ws.Range("E18").Value = "  +2.78%  "
$ws.Range("D19").Value = "61.872.94"
$ws.Range("E19").Value = "  +4.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.60"
$ws.Range("E20").Value = "  +4.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.94"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.45"
$ws.Range("E23").Value = "  +5.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.82"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.02"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +5.57%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  +4.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.17"
$ws.Range("E30").Value = "  +5.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.80"
$ws.Range("E31").Value = "  +10.81%  "
$ws.Range("E32").Value = "  +13.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.80"
$ws.Range("E33").Value = "  +3.97%  "
$ws.Range("E34").Value = "  +4.80%  "
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  +3.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.19"
$ws.Range("E37").Value = "  +5.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.38"
$ws.Range("E38").Value = "  +2.09%  "
$ws.Range("E39").Value = "  +9.36%  "
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "430.35"
$ws.Range("E41").Value = "  +6.79%  "
$ws.Range("E42").Value = "  +5.75%  "
$ws.Range("D43").Value = "2.850.00"
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.273"
$ws.Range("E44").Value = "  +7.46%  "
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.13"
$ws.Range("E46").Value = "  +6.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.69"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.07"

# Rows 47/48: Arweave and USDe swap places, with updated price/volume
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.26"
$ws.Range("E47").Value = "  +4.42%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.03%  "
